$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: "_old" columns -> "_FV2410", "_new" columns -> "_FV2504" ---
$headerNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J: old -> FV2410
for ($i = 0; $i -lt $headerNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($headerNames[$i])_FV2410"
}

# Column K stays "diff" (unchanged)

# Columns L-U: new -> FV2504
for ($i = 0; $i -lt $headerNames.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($headerNames[$i])_FV2504"
}

# --- 2. Freeze the top row (header row) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the used range into an Excel Table (ListObject) named "Table1" ---
# Stash the header row's existing direct formatting in a scratch row far away
# so it can be restored after the table is created, without Excel auto-generating
# a header dxf (which only happens when the header range already carries direct
# formatting at the moment the table is created).
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A70:U70")

$headerRange.Copy()
$scratchRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.Style = "Normal"

$dataRange = $ws.Range("A1:U68")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Restore the original header formatting now that the table exists
$scratchRange.Copy()
$headerRange.PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = $false

$scratchRange.Clear() | Out-Null

$ws.Range("A1").Select()
